$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: company_name id changes from "2" to "1", and its data columns get
# refreshed with Santam's newer dataset (same as row 3 below).
# Force text (not numeric) via the classic leading-apostrophe trick, then
# clear the resulting quote-prefix formatting so no stray style is left
# behind (matches the source, which keeps this cell style-less).
$ws.Range("B2").Value = "'1"
$ws.Range("B2").ClearFormats()

$ws.Range("D2").Value = 0.0669
$ws.Range("E2").Value = -0.0742
$ws.Range("G2").Value = 0.1657800863633519
$ws.Range("H2").Value = 0.1657800863633519
$ws.Range("I2").Value = 0.09856686901558295
$ws.Range("J2").Value = 0.05963709721951238
$ws.Range("K2").Value = 65.2
$ws.Range("L2").Value = 0.04080355466549847
$ws.Range("M2").Value = 80.8
$ws.Range("N2").Value = 0.04116358449233277
$ws.Range("O2").Value = 1.239263803680982
$ws.Range("P2").Value = 70.7
$ws.Range("Q2").Value = 0.03601813643079117
$ws.Range("R2").Value = 1.084355828220859
$ws.Range("S2").Value = 10.09999999999999
$ws.Range("T2").Value = 0.1249999999999999
$ws.Range("U2").Value = 323.5
$ws.Range("V2").Value = 0.164807173060268
$ws.Range("W2").Value = 0.1007727975270479
$ws.Range("X2").Value = 0.06327792465038341
$ws.Range("Y2").Value = 0.0374948728766645
$ws.Range("Z2").Value = 2.476212614287928
$ws.Range("AA2").Value = 0.1476741324144721
$ws.Range("AB2").Value = 0.06022275862014182
$ws.Range("AC2").Value = 0.08745137379433029
$ws.Range("AD2").Value = 219.8
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 219.8
$ws.Range("AG2").Value = -103.7
$ws.Range("AH2").Value = 0.1007009666926284
$ws.Range("AI2").Value = 0.2836495031616983
$ws.Range("AJ2").Value = -0.05577667814113597
$ws.Range("AK2").Value = -0.2297297297297297
$ws.Range("AL2").Value = 19.1
$ws.Range("AM2").Value = 19.1
$ws.Range("AN2").Value = 1.349294045426642
$ws.Range("AO2").Value = 8.246073298429319
$ws.Range("AP2").Value = -0.6365868631062001
$ws.Range("AQ2").Value = 8.246073298429319

# Row 3: company renamed from "Indequity Group Limited (JSE:IDQ)" to
# "Santam Ltd (JSE:SNT)" with its data columns updated to the same new
# dataset.
$ws.Range("B3").Value = "Santam Ltd (JSE:SNT)"

$ws.Range("D3").Value = 0.0669
$ws.Range("E3").Value = -0.0742
$ws.Range("G3").Value = 0.1657800863633519
$ws.Range("H3").Value = 0.1657800863633519
$ws.Range("I3").Value = 0.09856686901558295
$ws.Range("J3").Value = 0.05963709721951238
$ws.Range("K3").Value = 65.2
$ws.Range("L3").Value = 0.04080355466549847
$ws.Range("M3").Value = 80.8
$ws.Range("N3").Value = 0.04116358449233277
$ws.Range("O3").Value = 1.239263803680982
$ws.Range("P3").Value = 70.7
$ws.Range("Q3").Value = 0.03601813643079117
$ws.Range("R3").Value = 1.084355828220859
$ws.Range("S3").Value = 10.09999999999999
$ws.Range("T3").Value = 0.1249999999999999
$ws.Range("U3").Value = 323.5
$ws.Range("V3").Value = 0.164807173060268
$ws.Range("W3").Value = 0.1007727975270479
$ws.Range("X3").Value = 0.06327792465038341
$ws.Range("Y3").Value = 0.0374948728766645
$ws.Range("Z3").Value = 2.476212614287928
$ws.Range("AA3").Value = 0.1476741324144721
$ws.Range("AB3").Value = 0.06022275862014182
$ws.Range("AC3").Value = 0.08745137379433029
$ws.Range("AD3").Value = 219.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 219.8
$ws.Range("AG3").Value = -103.7
$ws.Range("AH3").Value = 0.1007009666926284
$ws.Range("AI3").Value = 0.2836495031616983
$ws.Range("AJ3").Value = -0.05577667814113597
$ws.Range("AK3").Value = -0.2297297297297297
$ws.Range("AL3").Value = 19.1
$ws.Range("AM3").Value = 19.1
$ws.Range("AN3").Value = 1.349294045426642
$ws.Range("AO3").Value = 8.246073298429319
$ws.Range("AP3").Value = -0.6365868631062001
$ws.Range("AQ3").Value = 8.246073298429319

# Row 4 (the old Santam Ltd row) is removed entirely - the sheet now only
# has 2 data rows, so delete the entire row 4 and shift rows up.
$ws.Rows("4").Delete()
